$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.397.43"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "3.840.23"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.13"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "3.839.23"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.65%  "
$ws.Range("E13").Value = "  +5.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.62%  "
$ws.Range("D15").Value = "4.483.84"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "3.849.79"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "68.453.47"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "470.34"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("E23").Value = "  -1.77%  "
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("E26").Value = "  -3.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "3.989.75"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.69"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("E35").Value = "  -2.52%  "
$ws.Range("D36").Value = "3.805.02"
$ws.Range("E36").Value = "  -1.31%  "
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("E38").Value = "  +10.05%  "
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("E44").Value = "  -4.68%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "414.41"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "47.12"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000291"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.31%  "
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.77%  "
